$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2021" column (O), mirroring the formatting
# of the existing "2020" column (N) for every populated row.
$rows = 3,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25

foreach ($r in $rows) {
    $src = $ws.Range("N" + $r)
    $dst = $ws.Range("O" + $r)
    $src.Copy($dst)
}

$ws.Range("O3").Value2 = 2021
$ws.Range("O5").Value2 = 2148.2
$ws.Range("O6").Value2 = 109.5
$ws.Range("O7").Value2 = 210.1
$ws.Range("O8").Value2 = 196
$ws.Range("O9").Value2 = 209
$ws.Range("O10").Value2 = 300.2
$ws.Range("O11").Value2 = 302.9
$ws.Range("O12").Value2 = 786
$ws.Range("O13").Value2 = 27.7
$ws.Range("O14").Value2 = 6.8
$ws.Range("O16").Value2 = 26.9
$ws.Range("O17").Value2 = 15.9
$ws.Range("O18").Value2 = 21.7
$ws.Range("O19").Value2 = 29.9
$ws.Range("O20").Value2 = 30.2
$ws.Range("O21").Value2 = 24
$ws.Range("O22").Value2 = 31.6
$ws.Range("O23").Value2 = 30.3
$ws.Range("O24").Value2 = 20.7
$ws.Range("O25").Value2 = 12

# Update the remembered selection to match the saved workbook state.
$ws.Range("Q20").Select()
